# Apply cryptos list price/volume update (commit: "Updated cryptos list on Sat Oct 28 11:35:15 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (matches source inlineStr cells) without leaving
# a stray NumberFormat-only style behind -- we restash the cells original Style
# object and reapply it after the write so the saved file keeps the same style index.
function Set-TextValue($cell, [string]$value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Row 2
$ws.Range("D2").Value = "34.277.20"
$ws.Range("E2").Value = "  +0.74%  "

# Row 3
$ws.Range("D3").Value = "1.794.39"
$ws.Range("E3").Value = "  +0.82%  "

# Row 4
$ws.Range("E4").Value = "  -0.22%  "

# Row 5
Set-TextValue $ws.Range("D5") "227.31"
$ws.Range("E5").Value = "  +0.77%  "

# Row 6
$ws.Range("E6").Value = "  +0.14%  "

# Row 7
$ws.Range("E7").Value = "  -0.21%  "

# Row 8
Set-TextValue $ws.Range("D8") "32.38"
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.294"
$ws.Range("E9").Value = "  +3.52%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.0695"
$ws.Range("E10").Value = "  -1.60%  "

# Row 11
$ws.Range("E11").Value = "  +0.61%  "

# Row 12
$ws.Range("D12").Value = "2.053.56"
$ws.Range("E12").Value = "  +0.81%  "

# Row 13
Set-TextValue $ws.Range("D13") "11.48"
$ws.Range("E13").Value = "  +5.07%  "

# Row 14
$ws.Range("D14").Value = "1.795.05"
$ws.Range("E14").Value = "  -0.72%  "

# Row 15
$ws.Range("D15").Value = "34.228.04"
$ws.Range("E15").Value = "  +0.74%  "

# Row 16
$ws.Range("E16").Value = "  +1.09%  "

# Row 17
$ws.Range("E17").Value = "  +2.17%  "

# Row 18
Set-TextValue $ws.Range("D18") "68.13"
$ws.Range("E18").Value = "  +0.78%  "

# Row 19
Set-TextValue $ws.Range("D19") "246.19"
$ws.Range("E19").Value = "  +1.21%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0784"
$ws.Range("E20").Value = "  +0.45%  "

# Row 21
$ws.Range("E21").Value = "  +2.67%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.998"
$ws.Range("E22").Value = "  -0.23%  "

# Row 23
$ws.Range("E23").Value = "  +1.43%  "

# Row 24
$ws.Range("E24").Value = "  -0.83%  "

# Row 25
Set-TextValue $ws.Range("D25") "161.74"
$ws.Range("E25").Value = "  +1.08%  "

# Row 26
Set-TextValue $ws.Range("D26") "7.19"
$ws.Range("E26").Value = "  +2.57%  "

# Row 27
Set-TextValue $ws.Range("D27") "16.35"
$ws.Range("E27").Value = "  +0.73%  "

# Row 28
$ws.Range("E28").Value = "  +1.83%  "

# Row 29
$ws.Range("E29").Value = "  -0.06%  "

# Row 30
$ws.Range("E30").Value = "  +1.29%  "

# Row 32
Set-TextValue $ws.Range("D32") "3.68"
$ws.Range("E32").Value = "  +1.95%  "

# Row 33
Set-TextValue $ws.Range("D33") "3.63"
$ws.Range("E33").Value = "  +3.96%  "

# Row 35
$ws.Range("D35").Value = "1.445.00"
$ws.Range("E35").Value = "  +4.17%  "

# Row 36
$ws.Range("E36").Value = "  +1.15%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.0193"
$ws.Range("E37").Value = "  +3.50%  "

# Row 38
$ws.Range("E38").Value = "  +10.07%  "

# Row 39
$ws.Range("E39").Value = "  -0.66%  "

# Row 40
Set-TextValue $ws.Range("D40") "80.95"
$ws.Range("E40").Value = "  +4.48%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.928"
$ws.Range("E41").Value = "  +2.23%  "

# Row 42
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("E43").Value = "  +0.35%  "

# Row 44
Set-TextValue $ws.Range("D44") "13.34"
$ws.Range("E44").Value = "  +7.26%  "

# Row 45
$ws.Range("D45").Value = "0.0₆0140"
$ws.Range("E45").Value = "  -0.38%  "

# Row 46
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D46") "0.0510"
$ws.Range("E46").Value = "  +2.79%  "

# Row 47
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D47") "6.07"
$ws.Range("E47").Value = "  +4.41%  "

# Row 48
$ws.Range("E48").Value = "  -1.08%  "

# Row 49
Set-TextValue $ws.Range("D49") "108.09"
$ws.Range("E49").Value = "  +0.69%  "

# Row 50
$ws.Range("D50").Value = "1.954.25"
$ws.Range("E50").Value = "  +0.89%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.998"
$ws.Range("E51").Value = "  -0.19%  "
